# Update the Metadata sheet (sheet1) and Include #0 sheet (sheet2)
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include #0")

# Metadata sheet updates
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-qualified-discrepancy"
$wsMeta.Range("B4").Value = "FRMedicationReconciliationQualifiedDiscrepancy"
$wsMeta.Range("B5").Value = "value set Interop'Santé - qualification de la divergence identifiée sur une ligne de traitement d'une FCT"
$wsMeta.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$wsMeta.Range("B11").Value = "FRANCE"

# Include #0 sheet updates
$wsInclude.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-discrepancy"
